# Bug fix: connector lines were saved with a near-zero width (1 EMU,
# i.e. <a:ln w="1">) instead of a sane hairline width of 1pt
# (12700 EMU == <a:ln w="12700">). Walk every shape on the slide and,
# for each straight-line connector, set its Line.Weight to 1 point so
# PowerPoint persists w="12700" in the XML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoConnector = 9

$updated = 0
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq $msoConnector) {
        $shp.Line.Weight = 1
        $updated = $updated + 1
    }
}

Write-Host "Connector lines updated:" $updated
